$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the previously used range first
$ws.Cells.Clear()

# Header row (order chosen to reproduce shared-string insertion order)
$ws.Range("D1").Value = "p1_cap"
$ws.Range("C1").Value = "l_power"
$ws.Range("B1").Value = "area"
$ws.Range("A1").Value = "macro"
$ws.Range("E1").Value = "p2_cap"

# Data rows: macro, area, l_power, p1_cap, p2_cap
$data = @(
    @("and",    1.0640000000000001, 25.07, 0.91810000000000003, 0.97460000000000002),
    @("or",     1.0640000000000001, 22.69, 0.94679999999999997, 0.94189999999999996),
    @("nand",   0.79800000000000004, 17.39, 1.599, 1.6641999999999999),
    @("nor",    0.79800000000000004, 21.2, 1.7144999999999999, 1.6513),
    @("not",    0.53200000000000003, 14.35, 1.7001999999999999, 0),
    @("xor",    1.5960000000000001, 36.159999999999997, 2.2321, 2.4115000000000002),
    @("xnor",   1.5960000000000001, 36.44, 2.2328000000000001, 2.5735999999999999),
    @("fflopd", 4.5220000000000002, 79.11, 0.94969999999999999, 1.1403000000000001),
    @("buf",    0.79800000000000004, 21.44, 0.97470000000000001, 0)
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $r++
}

$ws.Range("F14").Select()
